$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values: Day (date serial) + 24 hourly prices + daily avg + slot info
$ws.Range("A2").Value = 46020
$ws.Range("B2").Value = 99.58
$ws.Range("C2").Value = 91.33
$ws.Range("D2").Value = 89.76000000000001
$ws.Range("E2").Value = 88.03
$ws.Range("F2").Value = 89.98999999999999
$ws.Range("G2").Value = 92.06
$ws.Range("H2").Value = 101.29
$ws.Range("I2").Value = 116.99
$ws.Range("J2").Value = 127.17
$ws.Range("K2").Value = 114.06
$ws.Range("L2").Value = 102.99
$ws.Range("M2").Value = 93.27
$ws.Range("N2").Value = 89.73999999999999
$ws.Range("O2").Value = 89.59
$ws.Range("P2").Value = 95.16
$ws.Range("Q2").Value = 95.77
$ws.Range("R2").Value = 102.85
$ws.Range("S2").Value = 113.01
$ws.Range("T2").Value = 123.85
$ws.Range("U2").Value = 125.01
$ws.Range("V2").Value = 158.14
$ws.Range("W2").Value = 147.55
$ws.Range("X2").Value = 120.86
$ws.Range("Y2").Value = 104.73
$ws.Range("Z2").Value = 107.2

$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 132.82
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 152.84
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 124.43
$ws.Range("AG2").Value = "0h-23h"
